# 9.5.1.xlsx — add a 2023 data column (L) to the table on sheet 1.
#
# Source OOXML diff adds column L (year 2023) to rows 3-5 of the table,
# bumps every row's "spans" from 1:11 to 1:12 and the sheet <dimension> from
# A1:K13 to A1:L13, grows row 5's height to fit the (slightly) taller new
# header text, and drops the stale <selection activeCell="J12".../> that was
# left over from whoever last edited the sheet in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New column L: clone the formatting from column K (same row) then set
#     the new values, so the new cells inherit the exact same style indices
#     (borders/number formats/fonts) that K already uses on each row. ---

# Row 3 — empty bottom-border spacer cell under the title.
$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial(-4122)

# Row 4 — year header "2023".
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("L4").Value = 2023

# Row 5 — the data value for 2023.
$ws.Range("K5").Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("L5").Value = 0.11972285283622097

# Row 5 grew a bit taller (wrapped header text) in the authored version.
$ws.Rows.Item(5).RowHeight = 40.5

# Clear the stray J12 selection left in the saved view state; reset to A1
# (the workbook default) since that's the closest state to "no selection".
$ws.Range("A1").Select()
